$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column width fix for column AA (27): closest reachable value to target 12.7109375
$ws.Columns.Item(27).ColumnWidth = 11.833333333333334

# Cell value updates
$ws.Range("B1").Value = 0.99760234160146521
$ws.Range("F1").Value = 0.77607207723886007
$ws.Range("BI1").Value = 0.63099463950554058
$ws.Range("V2").Value = 0.77242586039908656
$ws.Range("AP2").Value = 0.83206319017012631
$ws.Range("AR2").Value = 0.95926988189581452
$ws.Range("BK2").Value = 0.88285931598273004
$ws.Range("H3").Value = 0.72489381406357412
$ws.Range("L3").Value = 0.9611926258137109
$ws.Range("V3").Value = 0.93272539818587097
$ws.Range("BP3").Value = 0.80725941913082133
$ws.Range("BD4").Value = 0.97131918194673683
$ws.Range("C5").Value = 0.92999031812534827
$ws.Range("D5").Value = 0.97565328816883379
$ws.Range("Z5").Value = 0.89218304147574334
$ws.Range("AN5").Value = 0.98388894430718787
$ws.Range("AG6").Value = 0.74141703468960796
$ws.Range("AF7").Value = 0.82859963162254802
$ws.Range("G9").Value = 0.96912974538947427
$ws.Range("H9").Value = 0.94090006657833059
$ws.Range("O9").Value = 0.88435413443829036
$ws.Range("AI9").Value = 0.93050050927153716
$ws.Range("BJ9").Value = 0.8121393418714975
$ws.Range("BE10").Value = 0.57282131745185727
$ws.Range("BL10").Value = 0.9747620600112944
$ws.Range("AI11").Value = 0.76419277467056523
$ws.Range("AT11").Value = 0.91655104765488216
$ws.Range("K13").Value = 0.58558490730272084
$ws.Range("Y13").Value = 0.74328592116531333
$ws.Range("M14").Value = 0.67712342234046485
$ws.Range("O14").Value = 0.64431465820992084
$ws.Range("M15").Value = 0.62622877288884271
$ws.Range("AN15").Value = 0.92635202019432761
$ws.Range("AO15").Value = 0.9447113394607467
$ws.Range("BI15").Value = 0.5868851175305192
$ws.Range("Z16").Value = 0.73621664492251626
$ws.Range("P18").Value = 0.9457659993721802
$ws.Range("Q18").Value = 0.68650013197796056
$ws.Range("Q19").Value = 0.81522961350378276
$ws.Range("R19").Value = 0.93641412115397493
$ws.Range("V19").Value = 0.82379727089629395
$ws.Range("AB19").Value = 0.91468451120493444
$ws.Range("BO19").Value = 0.75618989232814138
$ws.Range("AP20").Value = 0.65198771468076622
$ws.Range("AT20").Value = 0.74660557431321117
$ws.Range("AO21").Value = 0.95849287941554961
$ws.Range("AX21").Value = 0.69615479683753234
$ws.Range("AL23").Value = 0.73634661472316609
$ws.Range("AQ23").Value = 0.83219289246860817
$ws.Range("C24").Value = 0.94538002538787635
$ws.Range("AS24").Value = 0.94179091859807307
$ws.Range("BG24").Value = 0.99785962657924299
$ws.Range("BC25").Value = 0.95368146973409562
$ws.Range("I26").Value = 0.63954060123522505
$ws.Range("L29").Value = 0.91747583258037846
$ws.Range("R29").Value = 0.96679213352342142
$ws.Range("AJ29").Value = 0.97274594172220263
$ws.Range("V30").Value = 0.77399096721355032
$ws.Range("AJ30").Value = 0.87060904086894908
$ws.Range("AQ30").Value = 0.66374208498550824
$ws.Range("AZ30").Value = 0.92891366431158873
$ws.Range("BG30").Value = 0.83863110064438073
$ws.Range("K31").Value = 0.61924786614908989
$ws.Range("AU31").Value = 0.98124701116759916
$ws.Range("BF31").Value = 0.65712486236788936
$ws.Range("Q32").Value = 0.73397130018388834
$ws.Range("T33").Value = 0.9497950759926681
$ws.Range("AX33").Value = 0.62381034257728207
$ws.Range("AX34").Value = 0.98703514814841231
$ws.Range("BH34").Value = 0.93273660773013567
$ws.Range("BN34").Value = 0.83901028797922861
$ws.Range("J35").Value = 0.90462301371035903
$ws.Range("BC36").Value = 0.98193352162881475
$ws.Range("BO36").Value = 0.96450972514762823
$ws.Range("Z37").Value = 0.65373255755783433
$ws.Range("AJ37").Value = 0.98802303245895007
$ws.Range("AM37").Value = 0.70199829258581414
$ws.Range("BD37").Value = 0.8305553245047641
$ws.Range("R38").Value = 0.68572156312566979
$ws.Range("AK38").Value = 0.83023618433469626
$ws.Range("AS38").Value = 0.85833019134883992
$ws.Range("AM40").Value = 0.78170994521367598
$ws.Range("V41").Value = 0.684857203277907
$ws.Range("AM41").Value = 0.74819695504387207
$ws.Range("BD41").Value = 0.98234696295907642
$ws.Range("BG42").Value = 0.81856609022455151
$ws.Range("AR43").Value = 0.86166742180182965
$ws.Range("BH43").Value = 0.9088894693761792
$ws.Range("T44").Value = 0.90587202770923636
$ws.Range("L45").Value = 0.84456678032729005
$ws.Range("AA46").Value = 0.91211909320257079
$ws.Range("BH46").Value = 0.76933428115167324
$ws.Range("L47").Value = 0.74623363950806509
$ws.Range("S48").Value = 0.9537694111063284
$ws.Range("AW48").Value = 0.99745228836996613
$ws.Range("AM49").Value = 0.8173468971797343
$ws.Range("BP49").Value = 0.7882602498944552
$ws.Range("AA50").Value = 0.75744083570373888
$ws.Range("BF50").Value = 0.98673200783820203
$ws.Range("C51").Value = 0.706231116755762
$ws.Range("L51").Value = 0.96604610688377324
$ws.Range("N51").Value = 0.82396188080980504
$ws.Range("AH52").Value = 0.69674426106610854
$ws.Range("O53").Value = 0.83628703130179283
$ws.Range("AZ53").Value = 0.93111187682932317
$ws.Range("BL53").Value = 0.87162884462882517
$ws.Range("AD54").Value = 0.93959194615569297
$ws.Range("AU54").Value = 0.72379160013800681
$ws.Range("BA54").Value = 0.92254869853589438
$ws.Range("BE55").Value = 0.82642523409232216
$ws.Range("AB56").Value = 0.94267336337804819
$ws.Range("K57").Value = 0.80789480980143147
$ws.Range("Q57").Value = 0.83594212218574548
$ws.Range("BD57").Value = 0.92936048257912884
$ws.Range("BF57").Value = 0.96275700939826891
$ws.Range("AC58").Value = 0.99120958124043557
$ws.Range("B59").Value = 0.63997464987429575
$ws.Range("AB59").Value = 0.98370340975522541
$ws.Range("BL59").Value = 0.72766146887690164
$ws.Range("BA60").Value = 0.67344032576139523
$ws.Range("D61").Value = 0.97375119774146346
$ws.Range("AH61").Value = 0.86317246446396412
$ws.Range("AK61").Value = 0.68017135873679013
$ws.Range("AC62").Value = 0.87154922028353399
$ws.Range("AU62").Value = 0.97010809219839755
$ws.Range("M63").Value = 0.97219132602138614
$ws.Range("AQ63").Value = 0.94597936359639456
$ws.Range("BJ63").Value = 0.86606289937129799
$ws.Range("B65").Value = 0.98620252584750578
$ws.Range("X65").Value = 0.63137672264315714
$ws.Range("AG65").Value = 0.67128069194027717
$ws.Range("BE65").Value = 0.82082354660273493
$ws.Range("AF66").Value = 0.71948426752648764
$ws.Range("BP66").Value = 0.81315086203468989
$ws.Range("AB67").Value = 0.79315812036494082
$ws.Range("Q68").Value = 0.80077943942617491
